$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing date cell's formatting down onto the new rows so the
# new cells reuse the same style (numFmtId 14, "m/d/yyyy") instead of a
# freshly-minted custom number format.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Append the new rows of data (Date, Error Count)
$ws.Range("A10").Value = 45971
$ws.Range("B10").Value = 4

$ws.Range("A11").Value = 45973
$ws.Range("B11").Value = 11

$ws.Range("A12").Value = 45974
$ws.Range("B12").Value = 16

# Update the selection to match the post-edit state
$ws.Range("H14").Select()
